# Added player component to dg
#
# The canonical OOXML stores shape geometry in EMU (914400 EMU = 1 inch =
# 72 points) but the PowerPoint COM object model reads/writes Left/Top/
# Width/Height in points as a 32-bit float, so naive EMU/12700 round-trips
# can land 1 EMU off after the engine's internal f32 truncation. The point
# values below were solved so that they reproduce the exact target EMU
# after that conversion.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Resize/reposition the big background rounded rectangle ("Model").
# ---------------------------------------------------------------------
$bg = Get-ShapeById $s.Shapes 118
$bg.Left   = 95.08913045826772
$bg.Top    = 137.83843239685038
$bg.Width  = 600.9108582417323
$bg.Height = 341.9363861527559

# ---------------------------------------------------------------------
# 2. Shift the "Team" sub-diagram cluster down/right to make room for
#    the new player component (all shapes move by the same delta).
# ---------------------------------------------------------------------
$moves = @{
    77  = @(468.58409118818895, 420.69062812125986)
    82  = @(439.46543887086614, 354.8711090622047)
    87  = @(303.42442324881887, 419.777877815748)
    88  = @(575.5064392929133,  420.1262054724409)
    90  = @(450.3596038992126,  440.4215851031496)
    91  = @(418.35984811968507, 427.39503487007875)
    92  = @(437.88172918346453, 434.1363678527559)
    94  = @(284.4373322346457,  437.77156070314965)
    95  = @(525.4103088606299,  426.51252748503936)
    96  = @(541.7610779622047,  433.09864811732285)
    98  = @(485.81565859133855, 382.7933808267717)
    101 = @(496.46543887086614, 396.614089988189)
}

foreach ($id in $moves.Keys) {
    $shp = Get-ShapeById $s.Shapes $id
    $xy = $moves[$id]
    $shp.Left = $xy[0]
    $shp.Top  = $xy[1]
}

# ---------------------------------------------------------------------
# 3. Add the new "player" component: a rectangle labelled "UniqueTagList"
#    plus the two elbow connectors that wire it into the diagram.
# ---------------------------------------------------------------------

# 3a. New bentConnector3 arrow (clone of the unconnected "Elbow Connector
#     71" shape, which already carries the right line/arrow/style and no
#     stray connection-site bindings).
$srcA = Get-ShapeById $s.Shapes 96
$newA = $srcA.Duplicate().Item(1)
$newA.Name = "Elbow Connector 85"
$newA.Adjustments.Item(1) = 0.5
$newA.Left   = 547.9281311362205
$newA.Top    = 377.5326081251969
$newA.Width  = 84.49047474094489
$newA.Height = 0.0002363204724409449
$newA.Rotation = 270
$newA.HorizontalFlip = -1

# 3b. New rectangle "UniqueTagList" (clone of the existing UniqueTagList
#     shape so text formatting / style matches exactly).
$srcB = Get-ShapeById $s.Shapes 57
$newB = $srcB.Duplicate().Item(1)
$newB.Left   = 607.8559875519685
$newB.Top    = 383.43528749055116
$newB.Width  = 77.05448538897637
$newB.Height = 22.502598825196852

# 3c. New bentConnector2 arrow (clone of the matching "Elbow Connector 85"
#     shape that already has the right rotation/flip and no bindings).
$srcC = Get-ShapeById $s.Shapes 74
$newC = $srcC.Duplicate().Item(1)
$newC.Left   = 563.9314880629921
$newC.Top    = 349.8173980748032
$newC.Width  = 69.31763849527559
$newC.Height = 17.682913785826774
